# Apply cell-level updates to cryptos.xlsx per commit diff (Wed Dec 13 14:26:51 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '41.378.49'
$ws.Range("E2").Value = '  -1.06%  '

# Row 3
$ws.Range("D3").Value = '2.185.47'
$ws.Range("E3").Value = '  -1.42%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.39'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.42%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.74%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '67.02'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -4.94%  '

# Row 8
$ws.Range("E8").Value = '  +0.10%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.611'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.94%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.30'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.99%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.51'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.98%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0935'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.35%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.00'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.38%  '

# Row 14
$ws.Range("E14").Value = '  -1.15%  '

# Row 15
$ws.Range("D15").Value = '2.514.22'
$ws.Range("E15").Value = '  -1.30%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.45'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.91%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.853'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.78%  '

# Row 18
$ws.Range("D18").Value = '2.186.71'
$ws.Range("E18").Value = '  -1.74%  '

# Row 19
$ws.Range("D19").Value = '41.288.73'
$ws.Range("E19").Value = '  -1.03%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0949'
$ws.Range("E20").Value = '  -1.79%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.86'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.13%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.11'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.02%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.45'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.02%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.17%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.86'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.81%  '

# Row 26
$ws.Range("E26").Value = '  +0.07%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.24'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -7.89%  '

# Row 28
$ws.Range("E28").Value = '  -5.63%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.69'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.55%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.16'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.47%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.55'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.97%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.20'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.17%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0784'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.18%  '

# Row 34
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.120'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.68%  '

# Row 35
$ws.Range("E35").Value = '  +3.59%  '

# Row 36
$ws.Range("E36").Value = '  -1.14%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.17'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.07%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.81'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.09%  '

# Row 39
$ws.Range("E39").Value = '  -3.48%  '

# Row 40
$ws.Range("E40").Value = '  -0.32%  '

# Row 41
$ws.Range("E41").Value = '  -3.08%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.57'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -6.24%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.05'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.24%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.75'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.70%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.22'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.15%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.195'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.02%  '

# Row 47
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.100'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.17%  '

# Row 48
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.51'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.61%  '

# Row 49
$ws.Range("E49").Value = '  -0.52%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.15'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.83%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.85'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.67%  '
